# Update cryptos list prices / volume percentages (D2:E51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @("30.915.74", "  +2.65%  ")
    3  = @("1.904.40", "  +1.05%  ")
    4  = @("1.000", "  +0.09%  ")
    5  = @("246.56", "  +1.29%  ")
    6  = @("0.9996", "  +0.08%  ")
    7  = @("0.5006", "  +0.56%  ")
    8  = @("0.2995", "  +2.05%  ")
    9  = @("0.06866", "  +3.83%  ")
    10 = @("1.910.27", "  +1.42%  ")
    11 = @("17.47", "  +3.53%  ")
    12 = @("0.07341", "  +2.29%  ")
    13 = @("91.73", "  +6.90%  ")
    14 = @("5.123", "  +5.47%  ")
    15 = @("0.6827", "  +2.83%  ")
    16 = @("30.899.68", "  +2.65%  ")
    17 = @("0.000008082", "  +2.00%  ")
    18 = @("13.40", "  +4.72%  ")
    19 = @("1.001", "  +0.17%  ")
    20 = @("2.151.92", "  +1.44%  ")
    21 = @("1.001", "  +0.26%  ")
    22 = @("4.877", $null)
    23 = @("184.22", "  +36.82%  ")
    24 = @("6.122", "  +9.21%  ")
    25 = @("9.386", "  +2.85%  ")
    26 = @("154.04", "  +1.43%  ")
    27 = @("18.70", "  +11.39%  ")
    28 = @("1.958", "  +2.28%  ")
    29 = @("1.398", "  +1.23%  ")
    30 = @("4.396", "  +5.51%  ")
    31 = @("0.08998", "  +3.55%  ")
    32 = @("4.076", "  +3.27%  ")
    33 = @("0.05292", "  +5.75%  ")
    34 = @("0.7504", "  +6.13%  ")
    35 = @("1.144", "  +3.36%  ")
    36 = @("2.691", "  +1.59%  ")
    37 = @("0.01927", "  +16.75%  ")
    38 = @("2.727", "  +0.80%  ")
    39 = @("2.199", "  +0.37%  ")
    40 = @("0.9440", "  +0.91%  ")
    41 = @("0.4411", "  +5.31%  ")
    42 = @("106.57", "  +3.96%  ")
    43 = @("5.852", "  -1.84%  ")
    44 = @("1.0000", "  +0.10%  ")
    45 = @("7.798", "  +4.08%  ")
    46 = @("0.1361", "  +8.38%  ")
    47 = @("0.05856", "  +2.81%  ")
    48 = @("0.3939", "  +5.97%  ")
    49 = @("8.636", "  +4.84%  ")
    50 = @("33.53", "  +3.33%  ")
    51 = @("1.396", "  +3.81%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]

    if ($vals[1] -ne $null) {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $vals[1]
    }
}
